# Edit chapter 3.1.11 docx per commit diff:
#  - paragraph "BaseX aN521b_..." gains w:proofErr spellStart/spellEnd around
#    "BaseX" and gramStart/gramEnd around "dokumentbeskrivelse.xq", and the
#    filename run is split at the last underscore.
#  - paragraph "BaseX aN521_..." gets the same treatment.
#  - final "ANTALL ... er tomme og uten dokumenter, og er lagt til som
#    vedlegg." paragraph loses the trailing clause, ending in "." instead,
#    and a new 3-column table (RegistreringsID / SystemID / MappeID header
#    row) plus a trailing empty paragraph is appended after it.

$d = $word.ActiveDocument
$wordNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

# --- Paragraph 6: "BaseX aN521b_registreringer_uten_dokumentbeskrivelse.xq ..." ---
$p6 = $d.Paragraphs(6)
$p6xml = '<w:p xmlns:w="' + $wordNs + '">' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>BaseX</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:t>aN521b_registreringer_uten_</w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:t>dokumentbeskrivelse.xq</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:t xml:space="preserve"> for opptelling av registreringer uten dokumentbeskrivelse som ikke utgår.</w:t></w:r>' +
  '</w:p>'
[void]$p6.Range.InsertXML($p6xml)

# --- Paragraph 8: "BaseX aN521_registreringer_uten_dokumentbeskrivelse.xq ..." ---
$p8 = $d.Paragraphs(8)
$p8xml = '<w:p xmlns:w="' + $wordNs + '">' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>BaseX</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> aN521_registreringer_uten_</w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:t>dokumentbeskrivelse.xq</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:t xml:space="preserve"> for samme som over, men skriver ut hvilke registreringer det gjelder til fil. Denne kjøres om det er nødvendig.</w:t></w:r>' +
  '</w:p>'
[void]$p8.Range.InsertXML($p8xml)

# --- Final "ANTALL ... tomme og uten dokumenter, og er lagt til som vedlegg." paragraph ---
$lastIndex = $d.Paragraphs.Count
$pLast = $d.Paragraphs($lastIndex)
$pLastXml = '<w:p xmlns:w="' + $wordNs + '">' +
  '<w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t>ANTALL</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> registreringer er tomme og uten dokumenter</w:t></w:r>' +
  '<w:r><w:t>.</w:t></w:r>' +
  '</w:p>'
[void]$pLast.Range.InsertXML($pLastXml)

# --- Insert new table + trailing empty paragraph right after that paragraph ---
$pLast = $d.Paragraphs($d.Paragraphs.Count)
$insertPoint = $d.Range($pLast.Range.End, $pLast.Range.End)
$tableXml = '<w:tbl xmlns:w="' + $wordNs + '">' +
  '<w:tblPr>' +
    '<w:tblStyle w:val="Tabellrutenett"/>' +
    '<w:tblW w:w="0" w:type="auto"/>' +
    '<w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>' +
  '</w:tblPr>' +
  '<w:tblGrid>' +
    '<w:gridCol w:w="3020"/>' +
    '<w:gridCol w:w="3020"/>' +
    '<w:gridCol w:w="3021"/>' +
  '</w:tblGrid>' +
  '<w:tr>' +
    '<w:tc>' +
      '<w:tcPr>' +
        '<w:tcW w:w="3020" w:type="dxa"/>' +
        '<w:tcBorders>' +
          '<w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/>' +
          '<w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/>' +
          '<w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/>' +
          '<w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/>' +
        '</w:tcBorders>' +
        '<w:hideMark/>' +
      '</w:tcPr>' +
      '<w:p>' +
        '<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>RegistreringsID</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
      '</w:p>' +
    '</w:tc>' +
    '<w:tc>' +
      '<w:tcPr>' +
        '<w:tcW w:w="3020" w:type="dxa"/>' +
        '<w:tcBorders>' +
          '<w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/>' +
          '<w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/>' +
          '<w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/>' +
          '<w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/>' +
        '</w:tcBorders>' +
        '<w:hideMark/>' +
      '</w:tcPr>' +
      '<w:p>' +
        '<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>SystemID</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
      '</w:p>' +
    '</w:tc>' +
    '<w:tc>' +
      '<w:tcPr>' +
        '<w:tcW w:w="3021" w:type="dxa"/>' +
        '<w:tcBorders>' +
          '<w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/>' +
          '<w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/>' +
          '<w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/>' +
          '<w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/>' +
        '</w:tcBorders>' +
        '<w:hideMark/>' +
      '</w:tcPr>' +
      '<w:p>' +
        '<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>MappeID</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
      '</w:p>' +
    '</w:tc>' +
  '</w:tr>' +
  '</w:tbl>' +
  '<w:p/>'
[void]$insertPoint.InsertXML($tableXml)

Write-Output "Edit complete"
